$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.855.83"
$ws.Range("E2").Value = "  +1.61%  "
$ws.Range("D3").Value = "1.836.69"
$ws.Range("E3").Value = "  +1.44%  "
$ws.Range("E4").Value = "  +0.47%  "
$ws.Range("E5").Value = "  +0.41%  "
$ws.Range("D6").Value = "'308.44"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Value = "'0.4644"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +3.15%  "
$ws.Range("D8").Value = "'0.3614"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.84%  "
$ws.Range("D9").Value = "'0.07111"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Value = "'0.9087"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.23%  "
$ws.Range("B11").Value = "WrappedEther"
$ws.Range("C11").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D11").Value = "1.917.87"
$ws.Range("E11").Value = "  +5.76%  "
$ws.Range("B12").Value = "TRON"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D12").Value = "'0.07686"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.04%  "
$ws.Range("B13").Value = "Solana"
$ws.Range("C13").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D13").Value = "'19.47"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.69%  "
$ws.Range("D14").Value = "'5.262"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.17%  "
$ws.Range("D15").Value = "'6.367"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.94%  "
$ws.Range("D16").Value = "'87.66"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.30%  "
$ws.Range("D17").Value = "'1.011"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.61%  "
$ws.Range("D18").Value = "'0.000008566"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.68%  "
$ws.Range("E19").Value = "  +0.40%  "
$ws.Range("D20").Value = "26.907.01"
$ws.Range("E20").Value = "  +1.61%  "
$ws.Range("D21").Value = "'14.24"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.59%  "
$ws.Range("E22").Value = "  +0.73%  "
$ws.Range("D23").Value = "'10.61"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.98%  "
$ws.Range("D24").Value = "'1.937"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.93%  "
$ws.Range("D25").Value = "'152.16"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.68%  "
$ws.Range("D26").Value = "'18.14"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.97%  "
$ws.Range("D27").Value = "'2.016"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.99%  "
$ws.Range("D28").Value = "'113.73"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.49%  "
$ws.Range("D29").Value = "'4.859"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.51%  "
$ws.Range("D30").Value = "'0.08843"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.88%  "
$ws.Range("D31").Value = "'3.200"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.04%  "
$ws.Range("B32").Value = "ARBITRUM"
$ws.Range("C32").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D32").Value = "'1.166"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +5.23%  "
$ws.Range("D33").Value = "'0.7401"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.28%  "
$ws.Range("B34").Value = "RenderToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D34").Value = "'2.776"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.74%  "
$ws.Range("D35").Value = "'4.442"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.08%  "
$ws.Range("E36").Value = "  +1.04%  "
$ws.Range("B37").Value = "MXToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D37").Value = "'2.964"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.52%  "
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").Value = "'0.01928"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.09%  "
$ws.Range("D39").Value = "'0.05143"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.66%  "
$ws.Range("D40").Value = "'0.5153"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.66%  "
$ws.Range("D41").Value = "'6.875"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.67%  "
$ws.Range("D42").Value = "'0.1508"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.32%  "
$ws.Range("D43").Value = "'8.081"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.44%  "
$ws.Range("D44").Value = "'10.46"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +5.23%  "
$ws.Range("D45").Value = "'1.008"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Value = "'0.4662"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.70%  "
$ws.Range("D47").Value = "'99.86"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.02%  "
$ws.Range("D48").Value = "'1.593"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.01%  "
$ws.Range("D49").Value = "'0.06029"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.77%  "
$ws.Range("D50").Value = "'64.25"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.99%  "
$ws.Range("D51").Value = "'36.13"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.56%  "
